$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting of the last existing data row (row 14) down into the new row 15
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(15).RowHeight = 15.75

# Fill in the new record: date, start time, finish time, and computed hours
$ws.Range("A15").Value = 44857
$ws.Range("B15").Value = 0.70833333333333337
$ws.Range("C15").Value = 0.78333333333333333
$ws.Range("D15").Formula = "=C15-B15"

# Reflect the cursor position left behind after entering the new row
$ws.Range("C16").Select()

$wb.Save()
